$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (theta_threshold_range) - this shifts row 6 (pie_threshold_range) up to row 5
$ws.Rows.Item(5).Delete()

# Update values for the remaining parameter rows
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 10.4

$ws.Range("B3").Value = 5.6
$ws.Range("C3").Value = 9.2

$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 1.3

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

$ws.Range("C4").Select()
